$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("https://leetcode.com/problems/target-sum/?envType=list&envId=er2c1j13", "8c186c94199711eea88ae3300d621ca4", "arrays", "Target Sum", "medium", "leetcode"),
    @("https://leetcode.com/problems/maximum-sum-circular-subarray/?envType=list&envId=er2c1j13", "9e6b97e0199711eea88ae3300d621ca4", "arrays", "Maximum Sum Circular Subarray", "medium", "leetcode"),
    @("https://leetcode.com/problems/longest-substring-without-repeating-characters/?envType=list&envId=efdk2wsi", "bfcc797c199711eea88ae3300d621ca4", "strings", "Longest Substring Without Repeating Characters", "medium", "leetcode"),
    @("https://leetcode.com/problems/longest-palindromic-substring/?envType=list&envId=efdk2wsi", "d317bc62199711eea88ae3300d621ca4", "strings", "Longest Palindromic Substring", "medium", "leetcode"),
    @("https://leetcode.com/problems/decode-string/?envType=list&envId=efdk2wsi", "dfdd21da199711eea88ae3300d621ca4", "strings", "Decode String", "medium", "leetcode"),
    @("https://leetcode.com/problems/longest-common-prefix/?envType=list&envId=efdk2wsi", "072e0b1e199811eea88ae3300d621ca4", "strings", "Longest Common Prefix", "easy", "leetcode"),
    @("https://leetcode.com/problems/repeated-substring-pattern/?envType=list&envId=efdk2wsi", "36ff972c199811eea88ae3300d621ca4", "strings", "Repeated Substring Pattern", "easy", "leetcode"),
    @("https://leetcode.com/problems/sort-characters-by-frequency/?envType=list&envId=efdk2wsi", "4cb30a90199811eea88ae3300d621ca4", "strings", "Sort Characters By Frequency", "medium", "leetcode"),
    @("https://leetcode.com/problems/longest-duplicate-substring/?envType=list&envId=efdk2wsi", "6742bc0c199811eea88ae3300d621ca4", "strings", "Longest Duplicate Substring", "hard", "leetcode"),
    @("https://leetcode.com/problems/reformat-the-string/?envType=list&envId=efdk2wsi", "77897c90199811eea88ae3300d621ca4", "strings", "Reformat The String", "easy", "leetcode"),
    @("https://leetcode.com/problems/split-two-strings-to-make-palindrome/", "e71f12f4199811eea88ae3300d621ca4", "twoPointers", "Split Two Strings to Make Palindrome", "medium", "leetcode"),
    @("https://leetcode.com/problems/3sum/", "f5df63a2199811eea88ae3300d621ca4", "twoPointers", "3Sum", "medium", "leetcode"),
    @("https://leetcode.com/problems/shortest-unsorted-continuous-subarray/", "0759382e199911eea88ae3300d621ca4", "twoPointers", "Shortest Unsorted Continuous Subarray", "medium", "leetcode"),
    @("https://leetcode.com/problems/rotate-array/", "14bc1a04199911eea88ae3300d621ca4", "twoPointers", "Rotate Array", "medium", "leetcode"),
    @("https://leetcode.com/problems/next-permutation/", "374ca1d8199911eea88ae3300d621ca4", "twoPointers", "Next Permutation", "medium", "leetcode"),
    @("https://leetcode.com/problems/squares-of-a-sorted-array/", "535cf4c2199911eea88ae3300d621ca4", "twoPointers", "Squares of a Sorted Array", "easy", "leetcode")
)

$startRow = 9
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
